# Update Name of Algo
# Apply targeted numeric corrections to the result_data_RandomForest sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2"  = 16.4844
    "A3"  = -21.9634
    "A14" = -21.77429999999999
    "A21" = -20.25909999999999
    "A23" = -20.68189999999998
    "A25" = -21.84729999999999
    "E25" = 17.1672
    "A26" = -21.05529999999997
    "E27" = 16.75889999999999
    "A29" = -21.11749999999999
    "E31" = 16.52259999999999
    "E39" = 16.01439999999999
    "E48" = 17.5208
    "E51" = 17.3507
    "E52" = 17.26740000000001
    "A53" = -21.8801
    "E55" = 16.6873
    "E56" = 16.5006
    "A57" = -22.67930000000001
    "E57" = 16.7029
    "A59" = -22.22939999999999
    "A69" = -21.61039999999999
    "E73" = 17.30280000000001
    "A79" = -20.76120000000001
    "A83" = -21.84069999999999
    "E89" = 17.28170000000001
    "E90" = 16.50629999999999
    "A91" = -21.49640000000001
    "E92" = 18.91290000000002
    "A93" = -21.22609999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$wb.Save()
